$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($shape, $newName) {
    # Re-duplicate the shape's own Range before touching .Name - inline
    # shapes that live deep inside a header/footer story (i.e. not the
    # very first run of that story) can otherwise report a stale handle
    # when the Name property is set directly off the HeaderFooter.Range
    # collection, so re-seat through a fresh Range first.
    $rng = $shape.Range.Duplicate()
    $rng.InlineShapes.Item(1).Name = $newName
}

# First-page header logo: BTec_Logo-Orange (image1.jpg -> image2.jpg)
$hdrFirst = $sec.Headers.Item(2)
Rename-InlineLogo $hdrFirst.Range.InlineShapes.Item(1) "image2.jpg"

# Default footer logo: Pearson Edexcel logo (image2.png -> image1.png)
$ftrDefault = $sec.Footers.Item(1)
Rename-InlineLogo $ftrDefault.Range.InlineShapes.Item(1) "image1.png"

# First-page footer logo: Pearson Edexcel logo (image2.png -> image1.png)
$ftrFirst = $sec.Footers.Item(2)
Rename-InlineLogo $ftrFirst.Range.InlineShapes.Item(1) "image1.png"
